$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at the bottom of the data block (rows grow from 9 to 13 data rows)
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(10).Insert()

# Rewrite the full data block (rows 2-13) to match the updated quiz seed data
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "IT001"
$ws.Cells.Item(2, 3).Value = "DIEMDANH"
$ws.Cells.Item(2, 4).Value = 44927
$ws.Cells.Item(2, 5).Value = 44941.00013888889
$ws.Cells.Item(2, 6).Value = 2
$ws.Cells.Item(2, 7).Value = 1
$ws.Cells.Item(2, 8).Value = 0

$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "IT001"
$ws.Cells.Item(3, 3).Value = "BEGINER"
$ws.Cells.Item(3, 4).Value = 44927
$ws.Cells.Item(3, 5).Value = 44941.508472222224
$ws.Cells.Item(3, 6).Value = 30
$ws.Cells.Item(3, 7).Value = 2
$ws.Cells.Item(3, 8).Value = 0

$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "IT002"
$ws.Cells.Item(4, 3).Value = "DIEMDANH"
$ws.Cells.Item(4, 4).Value = 44941
$ws.Cells.Item(4, 5).Value = 44946.25708333333
$ws.Cells.Item(4, 6).Value = 2
$ws.Cells.Item(4, 7).Value = 1
$ws.Cells.Item(4, 8).Value = 0

$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "IT002"
$ws.Cells.Item(5, 3).Value = "BASIC"
$ws.Cells.Item(5, 4).Value = 44941
$ws.Cells.Item(5, 5).Value = 44946.25708333333
$ws.Cells.Item(5, 6).Value = 30
$ws.Cells.Item(5, 7).Value = 2
$ws.Cells.Item(5, 8).Value = 0

$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "IT003"
$ws.Cells.Item(6, 3).Value = "DIEMDANH"
$ws.Cells.Item(6, 4).Value = 44942
$ws.Cells.Item(6, 5).Value = 44946.229166666664
$ws.Cells.Item(6, 6).Value = 2
$ws.Cells.Item(6, 7).Value = 1
$ws.Cells.Item(6, 8).Value = 0

$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "IT003"
$ws.Cells.Item(7, 3).Value = "ADV"
$ws.Cells.Item(7, 4).Value = 44942
$ws.Cells.Item(7, 5).Value = 44946.229166666664
$ws.Cells.Item(7, 6).Value = 30
$ws.Cells.Item(7, 7).Value = 2
$ws.Cells.Item(7, 8).Value = 0

$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = "IT004"
$ws.Cells.Item(8, 3).Value = "DIEMDANH"
$ws.Cells.Item(8, 4).Value = 44950
$ws.Cells.Item(8, 5).Value = 44958.333333333336
$ws.Cells.Item(8, 6).Value = 2
$ws.Cells.Item(8, 7).Value = 1
$ws.Cells.Item(8, 8).Value = 0

$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = "IT004"
$ws.Cells.Item(9, 3).Value = "DATABASE"
$ws.Cells.Item(9, 4).Value = 44950
$ws.Cells.Item(9, 5).Value = 44958.375
$ws.Cells.Item(9, 6).Value = 30
$ws.Cells.Item(9, 7).Value = 2
$ws.Cells.Item(9, 8).Value = 0

$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "IT005"
$ws.Cells.Item(10, 3).Value = "DIEMDANH"
$ws.Cells.Item(10, 4).Value = 44962
$ws.Cells.Item(10, 5).Value = 44969.416666666664
$ws.Cells.Item(10, 6).Value = 2
$ws.Cells.Item(10, 7).Value = 1
$ws.Cells.Item(10, 8).Value = 0

$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = "IT005"
$ws.Cells.Item(11, 3).Value = "INTERN"
$ws.Cells.Item(11, 4).Value = 44962
$ws.Cells.Item(11, 5).Value = 44969.416666666664
$ws.Cells.Item(11, 6).Value = 30
$ws.Cells.Item(11, 7).Value = 2
$ws.Cells.Item(11, 8).Value = 0

$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = "IT006"
$ws.Cells.Item(12, 3).Value = "DIEMDANH"
$ws.Cells.Item(12, 4).Value = 44970
$ws.Cells.Item(12, 5).Value = 45052.458333333336
$ws.Cells.Item(12, 6).Value = 2
$ws.Cells.Item(12, 7).Value = 1
$ws.Cells.Item(12, 8).Value = 0

$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "IT006"
$ws.Cells.Item(13, 3).Value = "JUNIOR"
$ws.Cells.Item(13, 4).Value = 44970
$ws.Cells.Item(13, 5).Value = 45052.458333333336
$ws.Cells.Item(13, 6).Value = 30
$ws.Cells.Item(13, 7).Value = 2
$ws.Cells.Item(13, 8).Value = 0

# Update the active selection to match the saved view state
[void]$ws.Range("F12").Select()
